$d = $word.ActiveDocument

# The template placeholder "{{ piping_type }}" needs its variable name
# changed to "pipeline_construction_type". The original text is split
# across several runs ("piping", "_", "type") with proofErr markers
# (gramStart/gramEnd) wrapped around the first two pieces. A simple
# Find & Replace merges the runs into the found range and Word drops the
# now-orphaned grammar proofing marks automatically.
$d.Content.Find.Execute("piping_type", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "pipeline_construction_type", 2)
